$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.457.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.28%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.554.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.553.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "

# Row 11
$ws.Range("E11").Value = "  -0.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.021.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.95%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.356.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.542.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.07%  "

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.86%  "

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.04%  "

# Row 22
$ws.Range("E22").Value = "  -4.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "

# Row 24
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "

# Row 26
$ws.Range("E26").Value = "  -5.91%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.97%  "

# Row 28
$ws.Range("E28").Value = "  -4.81%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0929"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.68%  "

# Row 31
$ws.Range("E31").Value = "  -2.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "485.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("E34").Value = "  -3.04%  "

# Row 35
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("E36").Value = "  +6.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.97%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.98%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.40%  "

# Row 40
$ws.Range("E40").Value = "  -3.83%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("E42").Value = "  -5.12%  "

# Row 43
$ws.Range("E43").Value = "  -4.96%  "

# Row 45
$ws.Range("E45").Value = "  -3.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "145.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.31%  "

# Row 48
$ws.Range("E48").Value = "  -4.57%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.533"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.67%  "

# Row 50
$ws.Range("E50").Value = "  -6.76%  "

# Row 51
$ws.Range("E51").Value = "  -2.14%  "
